$d = $word.ActiveDocument

function Get-PlainText($range) {
    return $range.Text.Replace([char]13, "").Replace([char]7, "").Trim()
}

# ---------------------------------------------------------------------------
# 1) "Converted development environment ..." achievement bullet:
#    rewrite the sentence to add "our", "an", change "onboard time from 12
#    hours to 1 hour" -> "onboarding time by >90%".
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "Converted development environment to a Vagrant environment with automated installation script that reduced developer onboard time from 12 hours to 1 hour",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Converted our development environment to a Vagrant environment with an automated installation script that reduced developer onboarding time by >90%",
    2)
if (-not $found1) { throw "could not find the 'Converted development environment...' sentence" }

# ---------------------------------------------------------------------------
# 2) "Collaborated on infrastructure migration ..." bullet: drop the trailing
#    "with automated deployments using Ansible and Jenkins" clause - that
#    content becomes its own bullet below.
# ---------------------------------------------------------------------------
$found2 = $d.Content.Find.Execute(
    "Collaborated on infrastructure migration from a managed hosting platform to AWS with automated deployments using Ansible and Jenkins",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Collaborated on infrastructure migration from a managed hosting platform to AWS ",
    2)
if (-not $found2) { throw "could not find the 'Collaborated on infrastructure migration...' sentence" }

# ---------------------------------------------------------------------------
# 3) Insert a brand new achievement bullet right after the "Collaborated ..."
#    one, describing the Hubot/Slack automation tooling. It must inherit the
#    same Achievement / numbering / spacing formatting as its neighbours.
# ---------------------------------------------------------------------------
$collabPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ((Get-PlainText $p.Range) -eq "Collaborated on infrastructure migration from a managed hosting platform to AWS") {
        $collabPara = $p
    }
}
if ($collabPara -eq $null) { throw "could not locate the trimmed 'Collaborated ...AWS' paragraph" }

$collabPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ((Get-PlainText $p.Range) -eq "") {
        $prev = $p.Previous()
        if (($prev -ne $null) -and ((Get-PlainText $prev.Range) -eq "Collaborated on infrastructure migration from a managed hosting platform to AWS")) {
            $newPara = $p
        }
    }
}
if ($newPara -eq $null) { throw "could not locate the freshly inserted blank paragraph" }

$newPara.Range.InsertBefore("Built automation tools using Ansible and Jenkins with a Hubot interface and Slack integration") | Out-Null

# ---------------------------------------------------------------------------
# 4) The blank "Company Name One" paragraph that used to directly follow the
#    AWS/Jenkins bullet (and still holds the _GoBack bookmark) is restyled
#    into an Achievement paragraph (no bullet numbering), replacing its tab
#    stops / 360-twip indent with a 810-twip left indent and 276-auto line
#    spacing, while keeping its bold paragraph mark.
# ---------------------------------------------------------------------------
$spacerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if (("$($p.Style.NameLocal)" -like "*Company Name One*") -and ((Get-PlainText $p.Range) -eq "")) {
        $prev = $p.Previous()
        if (($prev -ne $null) -and ((Get-PlainText $prev.Range) -eq "Built automation tools using Ansible and Jenkins with a Hubot interface and Slack integration")) {
            $spacerPara = $p
        }
    }
}
if ($spacerPara -eq $null) { throw "could not locate the blank Company Name One spacer paragraph" }

$spacerPara.Format.Style = "Achievement"
$spacerPara.Format.LeftIndent = 40.5
$spacerPara.Format.FirstLineIndent = 0
$spacerPara.Format.LineSpacingRule = 5
$spacerPara.Format.LineSpacing = 13.8

# Re-apply the bold/size paragraph-mark formatting that setting .Format.Style
# above resets: briefly insert a placeholder character, paint it, then
# remove it again so only the paragraph-mark run-properties survive.
$spacerPara.Range.InsertBefore("X") | Out-Null
$markRange = $spacerPara.Range
$markRange.Font.Bold = $true
$markRange.Font.Size = 10
$markRange.Font.SizeBi = 9
$placeholder = $d.Range($markRange.Start, $markRange.Start + 1)
$placeholder.Delete() | Out-Null

"done"
